$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 1).Value = '亮仔'
$ws.Cells.Item(2, 2).Value = '嘭地一声'
$ws.Cells.Item(2, 3).Value = '回音山'
$ws.Cells.Item(2, 4).Value = '德鲁伊'

$ws.Cells.Item(3, 1).Value = '亮仔'
$ws.Cells.Item(3, 2).Value = '邪能肖战'
$ws.Cells.Item(3, 3).Value = '回音山'
$ws.Cells.Item(3, 4).Value = '恶魔猎手'

$ws.Cells.Item(4, 1).Value = '吴工'
$ws.Cells.Item(4, 2).Value = '体育老师'
$ws.Cells.Item(4, 3).Value = '通灵学院'
$ws.Cells.Item(4, 4).Value = '战士'

$ws.Cells.Item(5, 1).Value = '吴工'
$ws.Cells.Item(5, 2).Value = '邀月'
$ws.Cells.Item(5, 3).Value = '丽丽（四川）'
$ws.Cells.Item(5, 4).Value = '骑士'

$ws.Cells.Item(6, 1).Value = '吴工'
$ws.Cells.Item(6, 2).Value = '黑魔仙豹哥'
$ws.Cells.Item(6, 3).Value = '死亡之翼'
$ws.Cells.Item(6, 4).Value = '死亡骑士'

$ws.Cells.Item(7, 1).Value = '段总'
$ws.Cells.Item(7, 2).Value = '生锈的斩牛刀'
$ws.Cells.Item(7, 3).Value = '伊森利恩'
$ws.Cells.Item(7, 4).Value = '盗贼'

$ws.Cells.Item(8, 1).Value = '段总'
$ws.Cells.Item(8, 2).Value = '飞翔的潼瑜'
$ws.Cells.Item(8, 3).Value = '伊森利恩'
$ws.Cells.Item(8, 4).Value = '死亡骑士'

$ws.Cells.Item(9, 1).Value = '舒总'
$ws.Cells.Item(9, 2).Value = 'Fountine'
$ws.Cells.Item(9, 3).Value = '图拉扬'
$ws.Cells.Item(9, 4).Value = '法师'

$ws.Cells.Item(10, 1).Value = '舒总'
$ws.Cells.Item(10, 2).Value = '天灵浴血'
$ws.Cells.Item(10, 3).Value = '诺兹多姆'
$ws.Cells.Item(10, 4).Value = '死亡骑士'

$ws.Cells.Item(11, 1).Value = '统皇'
$ws.Cells.Item(11, 2).Value = '焦糖扁可颂'
$ws.Cells.Item(11, 3).Value = '斯坦索姆'
$ws.Cells.Item(11, 4).Value = '骑士'

$ws.Cells.Item(12, 1).Value = '统皇'
$ws.Cells.Item(12, 2).Value = '本间芽衣芓'
$ws.Cells.Item(12, 3).Value = '斯坦索姆'
$ws.Cells.Item(12, 4).Value = '战士'

$ws.Cells.Item(13, 1).Value = '巨奶'
$ws.Cells.Item(13, 2).Value = '傻瓜观测'
$ws.Cells.Item(13, 3).Value = '影之哀伤'
$ws.Cells.Item(13, 4).Value = '牧师'

$ws.Cells.Item(14, 1).Value = '捷教授'
$ws.Cells.Item(14, 2).Value = '四个自信'
$ws.Cells.Item(14, 3).Value = '回音山'
$ws.Cells.Item(14, 4).Value = '法师'

$ws.Cells.Item(15, 1).Value = '蔡圣'
$ws.Cells.Item(15, 2).Value = '莱恩弗尔特'
$ws.Cells.Item(15, 3).Value = '神圣之歌'
$ws.Cells.Item(15, 4).Value = '猎人'

$ws.Cells.Item(16, 1).Value = '蔡圣'
$ws.Cells.Item(16, 2).Value = '亚妮艾丝'
$ws.Cells.Item(16, 3).Value = '神圣之歌'
$ws.Cells.Item(16, 4).Value = '牧师'

$ws.Cells.Item(17, 1).Value = '蔡圣'
$ws.Cells.Item(17, 2).Value = '亚里欧斯'
$ws.Cells.Item(17, 3).Value = '神圣之歌'
$ws.Cells.Item(17, 4).Value = '恶魔猎手'

$ws.Cells.Item(18, 1).Value = '元神'
$ws.Cells.Item(18, 2).Value = '阿瘫'
$ws.Cells.Item(18, 3).Value = '霜之哀伤'
$ws.Cells.Item(18, 4).Value = '萨满'

$ws.Range("B24").Select()